$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Label" column (H) ---
# Copy the header formatting (bold font + border + center/top alignment)
# from the neighboring "Success %" header cell so the new header matches
# the existing header style, then overwrite its value/text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Label"

# Populate the new "Label" column: 0 for Control patients, 1 for MDD patients
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1, 0, 0, 0, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}

# --- Refit values (NCDEs refit to individual patients) ---
$ws.Range("D2").Value = 0.607959685570197
$ws.Range("E2").Value = 0.607959685570197

$ws.Range("D3").Value = 0.5203086362631032
$ws.Range("E3").Value = 0.5203086362631032

$ws.Range("D4").Value = 0.1989881160033442
$ws.Range("E4").Value = 0.1989881160033442

$ws.Range("D6").Value = 0.4264742770157975
$ws.Range("E6").Value = 0.4264742770157975

$ws.Range("D11").Value = 0.4309938993296114
$ws.Range("E11").Value = 0.5690061006703886
$ws.Range("F11").Value = 0.8202784657478333

$ws.Range("D16").Value = 0.4564388209386097
$ws.Range("E16").Value = 0.4564388209386097

$ws.Range("D21").Value = 0.7605294680110037
$ws.Range("E21").Value = 0.2394705319889963
$ws.Range("F21").Value = 0.7699308395385742
